$d = $word.ActiveDocument

$replacements = @(
    @('733×2=1466', '222×4=888'),
    @('134×8=1072', '123×2=246'),
    @('161×7=1127', '587×9=5283'),
    @('749×8=5992', '152×4=608'),
    @('280×8=2240', '482×6=2892'),
    @('501×5=2505', '536×2=1072'),
    @('393×5=1965', '212×6=1272'),
    @('564×4=2256', '325×2=650'),
    @('498×5=2490', '411×4=1644'),
    @('137×3=411', '246×6=1476'),
    @('335×2=670', '269×2=538'),
    @('243×9=2187', '728×4=2912'),
    @('233×5=1165', '960×2=1920'),
    @('502×3=1506', '157×9=1413'),
    @('164×6=984', '807×9=7263'),
    @('324×2=648', '912×3=2736'),
    @('602×5=3010', '411×5=2055'),
    @('529×3=1587', '529×6=3174'),
    @('238×2=476', '678×9=6102'),
    @('236×8=1888', '750×5=3750'),
    @('348×9=3132', '691×4=2764'),
    @('139×4=556', '740×2=1480'),
    @('582×9=5238', '822×8=6576'),
    @('657×4=2628', '601×4=2404'),
    @('283×7=1981', '324×5=1620'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done."
